# Rename model Zone to Region
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zone")
$ws.Name = "Region"

# Reflect the author's resulting view state: the renamed sheet becomes the
# active tab, with a fresh selection on it.
$ws.Activate() | Out-Null
$ws.Range("J21").Select() | Out-Null
